# Update the "http methods" notes sheet: refresh the HouseholdRestcontroller
# section and add a new FamilyMemberRestcontroller section with basic CRUD
# endpoints, per the commit "updated end points for basic crud on family member".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up stray "Done" markers in column I (rows 3-4) ---
$ws.Range("I3").Value = ""
$ws.Range("I4").Value = ""

# --- Row 5 used to be the "list householdTypes" (GET) row; repurpose it as
#     the "create a household" (POST) row ---
$ws.Range("A5").Value = "POST"
$ws.Range("B5").Value = "/household/create"
$ws.Range("F5").Value = "Create a new household"

# --- Clear out the old family-member rows (8, 9, 10); rebuilt below ---
$ws.Range("A8").Value = ""
$ws.Range("B8").Value = ""
$ws.Range("A9").Value = ""
$ws.Range("B9").Value = ""
$ws.Range("A10").Value = ""
$ws.Range("B10").Value = ""

# --- New section title above the household table ---
$ws.Range("A1").Value = "HouseholdRestcontroller"

# --- Row 6 used to duplicate the create endpoint; repurpose as "update" ---
$ws.Range("B6").Value = "/household/update"
$ws.Range("F6").Value = "Update a household"

# --- New section title for the family member table ---
$ws.Range("A9").Value = "FamilyMemberRestcontroller"

# --- Column headers for the new table ---
$ws.Range("A10").Value = "HTTP method"
$ws.Range("F10").Value = "CRUD action"

# --- New rows 11-15: basic CRUD for family members ---
$ws.Range("A11").Value = "GET"
$ws.Range("A12").Value = "GET"
$ws.Range("A13").Value = "POST"
$ws.Range("A14").Value = "PUT"
$ws.Range("A15").Value = "DELETE"

$ws.Range("B11").Value = "/household/family/list"
$ws.Range("B12").Value = "/household/family/list/{familyMemberId}"
$ws.Range("B13").Value = "/household/family/create"
$ws.Range("B14").Value = "/household/family/update"
$ws.Range("B15").Value = "/household/family/delete"

$ws.Range("F11").Value = "Read a list of family members"
$ws.Range("F12").Value = "Read a single family member"
$ws.Range("F13").Value = "Create a new family member"
$ws.Range("F14").Value = "Update a family member"
$ws.Range("F15").Value = "Delete a family member"

# --- Bold the two section titles and the two header rows ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A2").Font.Bold = $true
$ws.Range("F2").Font.Bold = $true
$ws.Range("A9").Font.Bold = $true
$ws.Range("A10").Font.Bold = $true
$ws.Range("F10").Font.Bold = $true

# --- Page setup / selection cosmetics ---
$ws.PageSetup.Orientation = 1
$ws.Range("A19").Select() | Out-Null
